# Update "想去人数" (interested-people count) for the two affected events
# on both the "展览" and "全部类型" worksheets.
#   F9:  58  -> 59
#   F10: 406 -> 407

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F9").Value = 59
    $ws.Range("F10").Value = 407
}
